$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    5597.0385765805,
    5427.810091879,
    5378.673491737501,
    5348.731399226001,
    5372.321170922,
    5496.759895401,
    6441.069913339,
    7184.668538918501,
    8465.353500000001,
    14333.5385,
    15915.495,
    15771.6165,
    15824.06,
    15963.6205,
    16132.936,
    16102.835,
    16787.6375,
    16403.816,
    16006.496,
    14046.071,
    11183.005,
    9240.505000000001,
    6857.7675,
    6096.6435
)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
